$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the old "invitaciones-boda" row (row 43). This shifts rows 44:63 up to 43:62.
$ws.Rows.Item(43).Delete()

# 2. The rows that used to be 55:63 (now 54:62) change their display_style
#    (column G) from "circle" to "card".
for ($r = 54; $r -le 62; $r++) {
    $ws.Cells.Item($r, 7).Value = "card"
}

# 3. The "papeleria-personal" row (now row 62) gets a new image_url that points
#    at the renamed "invitaciones_papeleria" media folder.
$ws.Cells.Item(62, 5).Value = "/media/subcategory_images/invitaciones_papeleria/papeleria-personal.jpg"

# 4. Append the brand new "bodas" row at the bottom (row 63).
$ws.Cells.Item(63, 1).Value = "bodas"
$ws.Cells.Item(63, 2).Value = "Todo para tu boda"
$ws.Cells.Item(63, 3).Value = "invitaciones-papeleria"
$ws.Cells.Item(63, 4).Value = "Todo para tu boda"
$ws.Cells.Item(63, 5).Value = "/media/subcategory_images/invitaciones_papeleria/invitaciones-boda.jpg"
$ws.Cells.Item(63, 6).Value = 12
$ws.Cells.Item(63, 7).Value = "card"

# 5. Shrink the hidden _xlnm._FilterDatabase defined name by one row (59 -> 58)
#    to match the net removal of one data row.
$wb.Names.Item("_xlnm._FilterDatabase").RefersTo = "=subcategories_complete!`$A`$1:`$G`$58"

# 6. Update the saved view/selection state.
$ws.Application.ActiveWindow.ScrollRow = 33
$ws.Range("A58:A63").Select()
